# suppression du champ Batiment
# - grows the "Adresse" table cell by two extra blank lines and resizes/
#   repositions the table accordingly
# - enlarges / repositions the "nomBatiment" title textbox and bumps its
#   font size
# - removes the now redundant "description" textbox
# - shifts the adresse / dateDeConstruction / surfaceTotaleChauffe /
#   dateDeRenovation textboxes down (and adresse left) to fill the gap

# PowerPoint's Shape.Left/Top/Width/Height are Single-precision point
# values, while the underlying OOXML stores EMUs (1 pt = 12700 EMU).
# Converting naively (emu/12700.0) occasionally lands one EMU short once
# it is round-tripped through Single precision, so probe nearby Single
# values until we find one that reproduces the exact EMU on save.
function EmuToPt([double]$targetEmu) {
    $pts = $targetEmu / 12700.0
    for ($i = 0; $i -lt 400; $i++) {
        $f = [single]$pts
        $emu = [math]::Floor([double]$f * 12700.0)
        if ($emu -eq $targetEmu) { return $f }
        if ($emu -lt $targetEmu) { $pts += 0.0000005 } else { $pts -= 0.0000005 }
    }
    return [single]$pts
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Locate the shapes we need by name -------------------------------
$table = $null
$nomBatiment = $null
$description = $null
$adresse = $null
$dateDeConstruction = $null
$surfaceTotaleChauffe = $null
$dateDeRenovation = $null

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    switch ($shp.Name) {
        "Tableau 1"             { $table = $shp }
        "nomBatiment"            { $nomBatiment = $shp }
        "description"            { $description = $shp }
        "adresse"                { $adresse = $shp }
        "dateDeConstruction"     { $dateDeConstruction = $shp }
        "surfaceTotaleChauffe"   { $surfaceTotaleChauffe = $shp }
        "dateDeRenovation"       { $dateDeRenovation = $shp }
    }
}

# --- 1. Table: grow it downward and add two blank lines under Adresse -
# Add the two extra blank paragraphs first: the table's autofit height
# (the engine's `pptx_autofit_height_emu`-style relayout) then already
# lands on the correct new frame height as soon as we touch its
# position, so Height doesn't need to be (and must not be) forced
# separately - doing so would make the engine re-split the row heights
# proportionally instead of leaving the untouched rows alone.
$addrCell = $table.Table.Cell(1, 1)
$addrRange = $addrCell.Shape.TextFrame.TextRange
$addrRange.Text = "Adresse :" + [char]13 + [char]13 + [char]13 + [char]13

$table.Top = EmuToPt(1710425)

# --- 2. nomBatiment: widen/heighten, shift left, bigger font ----------
$nomBatiment.Left = EmuToPt(6957634)
$nomBatiment.Width = EmuToPt(5089783)
$nomBatiment.Height = EmuToPt(461665)
$nomBatiment.TextFrame.TextRange.Font.Size = 24

# --- 3. description: remove this textbox entirely ---------------------
$description.Delete()

# --- 4. adresse: move left + down --------------------------------------
$adresse.Left = EmuToPt(9048137)
$adresse.Top = EmuToPt(1663213)

# --- 5. dateDeConstruction: move down -----------------------------------
$dateDeConstruction.Top = EmuToPt(3543591)

# --- 6. surfaceTotaleChauffe: move down ---------------------------------
$surfaceTotaleChauffe.Top = EmuToPt(3130566)

# --- 7. dateDeRenovation: move down -------------------------------------
$dateDeRenovation.Top = EmuToPt(3939605)
